$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 8 data rows (years 2000,2002,2005,2007,2010,2012,2015,2017)
# in rows 2-9 below the header row 1. The update drops the first four years
# (2000,2002,2005,2007 -> old rows 2-5), keeping 2010/2012/2015/2017 which
# shift up to rows 2-5, and appends a new 2020 row at row 6.

# 1) Remove the obsolete year rows (2000, 2002, 2005, 2007).
$ws.Rows("2:5").Delete() | Out-Null

# After the delete, rows are: 2=2010, 3=2012, 4=2015, 5=2017 (unchanged values).
# 2) Copy the formatting of the last existing row (2017, row 5) down to the
#    new row 6 so the new row matches the sheet's existing look (bold/bordered
#    label cell, etc.).
$ws.Range("A5:S5").Copy() | Out-Null
$ws.Range("A6:S6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Populate the new 2020 row's values. The remaining columns (B, E, F, G,
#    H, J, K, L, M, N, O, P, Q, R, S) have no reported figure for 2020, so
#    they are left blank, matching the blank pattern used throughout the
#    rest of the sheet.
$ws.Cells.Item(6, 1).Value = "2020年"
$ws.Cells.Item(6, 3).Value = 1482404205.74338
$ws.Cells.Item(6, 4).Value = 19639627.5808391
$ws.Cells.Item(6, 9).Value = 1736253794.41957
